$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.887.51"
$ws.Range("E2").Value = "  -1.23%  "

# Row 3
$ws.Range("D3").Value = "2.194.99"
$ws.Range("E3").Value = "  -2.33%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.30%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "295.24"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -4.00%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "88.52"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -6.32%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.567"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.61%  "

# Row 8
$ws.Range("E8").Value = "  -0.10%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.480"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -8.78%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.13"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -7.66%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0769"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -5.37%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.103"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.40%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.75"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -5.83%  "

# Row 14
$ws.Range("D14").Value = "2.525.19"

# Row 15
$ws.Range("D15").Value = "2.255.74"
$ws.Range("E15").Value = "  -4.42%  "

# Row 16
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "12.98"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -5.02%  "

# Row 17
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.769"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -8.51%  "

# Row 18
$ws.Range("D18").Value = "43.490.09"
$ws.Range("E18").Value = "  -1.29%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0885"
$ws.Range("E19").Value = "  -8.11%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.81"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -9.19%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.66"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -14.81%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "62.84"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.65%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "230.10"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.11%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.75"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -12.67%  "

# Row 25
$ws.Range("E25").Value = "  +0.54%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.81"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -9.83%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.15"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.77%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "35.70"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -7.51%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.16"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -7.07%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.08"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -5.17%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "146.68"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -4.95%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.29"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -11.36%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.49"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -6.10%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0731"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -8.66%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.116"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.23%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.85"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -8.45%  "

# Row 37
$ws.Range("E37").Value = "  -6.80%  "

# Row 38
$ws.Range("E38").Value = "  -10.25%  "

# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0281"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -7.94%  "

# Row 40
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.50"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -8.28%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.07"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -11.36%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").ClearFormats()

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.94"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -11.31%  "

# Row 44
$ws.Range("D44").Value = "1.771.10"
$ws.Range("E44").Value = "  +1.53%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.61"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.16%  "

# Row 46
$ws.Range("B46").Value = "BitcoinSV"
$ws.Range("C46").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "72.64"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -10.14%  "

# Row 47
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.172"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -11.04%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "13.68"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +6.14%  "

# Row 49
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "91.50"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -8.14%  "

# Row 50
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.409.02"
$ws.Range("E50").Value = "  -2.37%  "

# Row 51
$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.47"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -9.27%  "
